$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RequirementInfo")
$ws.Name = "RequirementInfoData"
$ws.Range("E11").Select()
